# Apply the "Program" sheet data edits:
#  - B3: Devops -> corejava
#  - C3: BDD -> OOPs concepts
#  - D4: Devops -> corejava
#  - B5: DelTestOne -> Cybersecurity
#  - Selection moves to D4
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program")

# Order matters for how new entries land in the shared-strings table:
# C3 first so "OOPs concepts" is interned before "corejava".
$ws.Range("C3").Value = "OOPs concepts"
$ws.Range("B3").Value = "corejava"
$ws.Range("D4").Value = "corejava"
$ws.Range("B5").Value = "Cybersecurity"

# B5 picks up the same "Arial 18" font styling already used by the other
# s="6" cells in this block (B3/B4/D4/E4) instead of the plain body style.
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)  # xlPasteFormats: formatting only, keep B5's value

$ws.Activate()
$ws.Range("D4").Select()
